$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted as row 10, pushing every existing
# record from row 10 downward down by one row (old row 10 -> new row 11,
# ..., old row 35 -> new row 36).
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Macroferia Regional de Talca"
$ws.Range("C10").Value = "Maule"
$ws.Range("D10").Value = 44742
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 100112043
$ws.Range("G10").Value = "Pepino dulce"
$ws.Range("H10").Value = "Cultivar IV Región"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 14000
$ws.Range("N10").Value = "$/bandeja 18 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 778
$ws.Range("Q10").Value = 18
$ws.Range("R10").Value = "Hortaliza"
